# Add a new "2022-Q4" quarter sheet (with its fund holdings) right after
# the "总计" (summary) sheet, and update the summary sheet with the new
# quarter's totals, pushing all older quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计" (the first tab).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet    = $wb.Worksheets.Item("2022-Q3")

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Header row (columns B..H), matching the other quarterly sheets.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (A = sequential 0-based index, B/D/E/F/G stored as text like
# the other sheets -- B keeps leading zeros on fund codes -- H stored as a
# number).
$newSheet.Range("B2:B4").NumberFormat = "@"
$newSheet.Range("D2:G4").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "513090"
$newSheet.Range("C2").Value = "易方达中证香港证券投资主题ETF"
$newSheet.Range("D2").Value = "11.28"
$newSheet.Range("E2").Value = "97.07"
$newSheet.Range("F2").Value = "7.73"
$newSheet.Range("G2").Value = "0.8719"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "004497"
$newSheet.Range("C3").Value = "前海开源多元策略灵活配置混合C"
$newSheet.Range("D3").Value = "1.79"
$newSheet.Range("E3").Value = "79.66"
$newSheet.Range("F3").Value = "3.62"
$newSheet.Range("G3").Value = "0.0648"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "004496"
$newSheet.Range("C4").Value = "前海开源多元策略灵活配置混合A"
$newSheet.Range("D4").Value = "1.30"
$newSheet.Range("E4").Value = "79.66"
$newSheet.Range("F4").Value = "3.62"
$newSheet.Range("G4").Value = "0.0471"
$newSheet.Range("H4").Value = 9

# Formatting: bold, centered/top-aligned, thin-bordered header row + index
# column, matching the look of the other quarterly sheets.
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexRange = $newSheet.Range("A2:A4")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: prepend the 2022-Q4 totals and
#    shift every older quarter down by one row (2020-Q4 becomes row 10).
# ---------------------------------------------------------------------
$quarters = @(
    @{ Label = "2022-Q4"; Count = 3; Value = 0.98 },
    @{ Label = "2022-Q3"; Count = 6; Value = 0.99 },
    @{ Label = "2022-Q2"; Count = 3; Value = 1.22 },
    @{ Label = "2022-Q1"; Count = 4; Value = 1.23 },
    @{ Label = "2021-Q4"; Count = 5; Value = 1.1 },
    @{ Label = "2021-Q3"; Count = 3; Value = 1.12 },
    @{ Label = "2021-Q2"; Count = 6; Value = 1.8 },
    @{ Label = "2021-Q1"; Count = 4; Value = 1.89 },
    @{ Label = "2020-Q4"; Count = 7; Value = 1.92 }
)

for ($i = 0; $i -lt $quarters.Count; $i++) {
    $row = 2 + $i
    $totalSheet.Range("A$row").Value = $i
    $totalSheet.Range("B$row").Value = $quarters[$i].Label
    $totalSheet.Range("C$row").Value = $quarters[$i].Count
    $totalSheet.Range("D$row").Value = $quarters[$i].Value
}

# The brand-new last row (A10) needs the same bold/centered/bordered look
# as the rest of column A on this sheet.
$lastIndexCell = $totalSheet.Range("A" + (1 + $quarters.Count))
$lastIndexCell.Font.Bold = $true
$lastIndexCell.HorizontalAlignment = -4108
$lastIndexCell.VerticalAlignment = -4160
$lastIndexCell.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3. Keep the originally-selected tab ("2020-Q4", the last sheet) active.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
